$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI script was re-run with updated TPM data. The sending
# cluster for this Slitrk2-Ptprs edge table is now "MuSCs" (was "ECs"), and
# the edge-weight / specificity metrics for every target cluster row were
# recomputed accordingly. Target cluster labels (column D) are unchanged.

# Sending cluster (column A) changes from "ECs" to "MuSCs" for all data rows
$ws.Range("A2:A7").Value = "MuSCs"

# Updated numeric values (G, H, M, N, O, P, Q, R, S, T) per row, recomputed with new TPM data

# Row 2 (ECs)
$ws.Range("G2").Value = 0.01213966666666667
$ws.Range("H2").Value = 0.036419
$ws.Range("M2").Value = 4.260872666666667
$ws.Range("N2").Value = 12.782618
$ws.Range("O2").Value = 0.09064705929364961
$ws.Range("P2").Value = 0.09064705929364959
$ws.Range("Q2").Value = 0.05172557388244445
$ws.Range("R2").Value = 0.4655301649420001
$ws.Range("S2").Value = 0.09064705929364961
$ws.Range("T2").Value = 0.09064705929364959

# Row 3 (FAPs)
$ws.Range("G3").Value = 0.01213966666666667
$ws.Range("H3").Value = 0.036419
$ws.Range("M3").Value = 20.524797
$ws.Range("N3").Value = 61.574391
$ws.Range("O3").Value = 0.4366505728284585
$ws.Range("P3").Value = 0.4366505728284584
$ws.Range("Q3").Value = 0.249164193981
$ws.Range("R3").Value = 2.242477745829
$ws.Range("S3").Value = 0.4366505728284585
$ws.Range("T3").Value = 0.4366505728284584

# Row 4 (Inflammatory-Mac)
$ws.Range("G4").Value = 0.01213966666666667
$ws.Range("H4").Value = 0.036419
$ws.Range("M4").Value = 8.931090666666666
$ws.Range("N4").Value = 26.793272
$ws.Range("O4").Value = 0.190002651698962
$ws.Range("P4").Value = 0.1900026516989619
$ws.Range("Q4").Value = 0.1084204636631111
$ws.Range("R4").Value = 0.9757841729679999
$ws.Range("S4").Value = 0.190002651698962
$ws.Range("T4").Value = 0.1900026516989619

# Row 5 (MuSCs)
$ws.Range("G5").Value = 0.01213966666666667
$ws.Range("H5").Value = 0.036419
$ws.Range("M5").Value = 6.457974333333333
$ws.Range("N5").Value = 19.373923
$ws.Range("O5").Value = 0.1373888468646722
$ws.Range("P5").Value = 0.1373888468646721
$ws.Range("Q5").Value = 0.07839765574855556
$ws.Range("R5").Value = 0.705578901737
$ws.Range("S5").Value = 0.1373888468646722
$ws.Range("T5").Value = 0.1373888468646721

# Row 6 (Neutrophils)
$ws.Range("G6").Value = 0.01213966666666667
$ws.Range("H6").Value = 0.036419
$ws.Range("M6").Value = 1.948535
$ws.Range("N6").Value = 5.845605
$ws.Range("O6").Value = 0.04145370713904261
$ws.Range("P6").Value = 0.0414537071390426
$ws.Range("Q6").Value = 0.02365456538833333
$ws.Range("R6").Value = 0.212891088495
$ws.Range("S6").Value = 0.04145370713904261
$ws.Range("T6").Value = 0.0414537071390426

# Row 7 (Resolving-Mac)
$ws.Range("G7").Value = 0.01213966666666667
$ws.Range("H7").Value = 0.036419
$ws.Range("M7").Value = 4.881814666666666
$ws.Range("N7").Value = 14.645444
$ws.Range("O7").Value = 0.1038571621752152
$ws.Range("P7").Value = 0.1038571621752152
$ws.Range("Q7").Value = 0.05926360278177777
$ws.Range("R7").Value = 0.5333724250359999
$ws.Range("S7").Value = 0.1038571621752152
$ws.Range("T7").Value = 0.1038571621752152
